$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("D2").Value = '27.314.86'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '1.828.03'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.89'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9988'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4490'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.77%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3790'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07528'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8877'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.07'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.65%  '
$ws.Range("D12").Value = '1.814.52'
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.778'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.80'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.411'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07124'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9994'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008821'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9984'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.17'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("D21").Value = '27.337.31'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.276'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").Value = '2.049.98'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.991'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.468'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.23'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.63'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.394'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.39'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08850'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7749'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.195'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.587'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9985'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.112'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01997'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05325'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.437'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5357'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1735'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.860'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.258'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +14.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.816'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5127'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.83'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '106.87'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.708'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9983'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06386'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.99%  '
